$d = $word.ActiveDocument

# This document had a "featured" screenshot inserted right under the title
# (the first inline picture in the document) plus the *same* screenshot
# pasted again as a whole extra paragraph right after the "Introducción"
# heading, and a third time right after the "Presentación del proyecto"
# heading. Those two extra copies are duplicated captures and must be
# removed completely (the image run, its own paragraph, and the paragraph
# mark), while the featured image at the top is left untouched.
#
# Walk the inline pictures from the end of the document towards the start
# (so deleting one doesn't shift the index of the ones still to process)
# and drop every picture-only paragraph except the very first picture in
# the document.
$shapeCount = $d.InlineShapes.Count
for ($i = $shapeCount; $i -ge 2; $i--) {
    $shape = $d.InlineShapes.Item($i)
    $shapeParagraph = $shape.Range.Paragraphs.Item(1)

    # Safety check: only remove the paragraph if it contains nothing but
    # the picture itself (mirrors the paragraphs removed in the diff,
    # which held a single centered <w:drawing> run and nothing else).
    if ($shapeParagraph.Range.InlineShapes.Count -eq 1 -and $shapeParagraph.Range.Text.Length -eq 1) {
        $shapeParagraph.Range.Delete()
    }
}
